$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Column C (AccessLevel) moves from text labels to numeric codes.
$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 4
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 3

# Selection moves to C4.
$ws.Range("C4").Select()
